$wb = $excel.ActiveWorkbook

# Rename sheets
$ws1 = $wb.Worksheets.Item("LIT")
$ws2 = $wb.Worksheets.Item("GRAMMAR")
$ws1.Name = "Lit 3 2019-2020"
$ws2.Name = "Grammar 4 2019-2020"

# Sheet1 ("Lit 3 2019-2020"): drop the Retake column (C) and the Year column (now D
# after the first delete, originally E), fill in student marks
$ws1.Range("C1").EntireColumn.Delete()
$ws1.Range("D1").EntireColumn.Delete()

$ws1.Range("A2:B6").NumberFormat = "@"

$ws1.Range("A2").Value = "Isa"
$ws1.Range("B2").Value = "Garvi"
$ws1.Range("C2").Value = 9

$ws1.Range("A3").Value = "Cris"
$ws1.Range("B3").Value = "Rodero"
$ws1.Range("C3").Value = 8

$ws1.Range("A4").Value = "Mariana"
$ws1.Range("B4").Value = "Santos"
$ws1.Range("C4").Value = 3

$ws1.Range("A5").Value = "Alejandra"
$ws1.Range("B5").Value = "Carretero"
$ws1.Range("C5").Value = 2

$ws1.Range("A6").Value = "Nela"
$ws1.Range("B6").Value = "Alberola"
$ws1.Range("C6").Value = 7
